$wb = $excel.ActiveWorkbook

# --- Sheet "goods init" ---
$goods = $wb.Worksheets.Item("goods init")

# Update existing quantities
$goods.Range("B2").Value = 13.6   # flour
$goods.Range("B3").Value = 10     # eggs
$goods.Range("B4").Value = 2.65   # butter
$goods.Range("B5").Value = 4.7    # cheese

# Add new row for sugar
$goods.Range("A9").Value = "sugar"
$goods.Range("B9").Value = 2
$goods.Range("C9").Value = "kg"

# --- Sheet "product init" ---
$products = $wb.Worksheets.Item("product init")

# Update existing quantities
$products.Range("C3").Value = 13  # bread
$products.Range("C8").Value = 11  # cheese pie

# Add new row for oreo
$products.Range("A9").Value = "oreo"
$products.Range("B9").Value = 2
$products.Range("C9").Value = 0
